# edit.ps1 - applies the "Add reference gradients and fix logo on all slides" change
# described by the target diff:
#   1. Slide 3 ("Анализ и рекомендации") background Rectangle 1: replace the flat
#      white a:solidFill with a two-stop linear a:gradFill (F8F9FA -> E8F5E8).
#   2. Slide 3 narrative textbox: rewrite the summary / bottleneck / next-step
#      paragraphs with the new analysis text.
#   3. Slide 9 ("Выводы и рекомендации") narrative textbox: rewrite the summary /
#      bottleneck / next-step paragraphs with the new analysis text.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Slide 3 background rectangle: solid white fill -> reference gradient fill
# ---------------------------------------------------------------------------
$slide3 = $p.Slides.Item(3)
$bgShape = $slide3.Shapes.Item(1)

# msoGradientHorizontal (2), variant 1 -> linear gradient at angle 0, which is
# the schema default for <a:lin ang>, matching the authored "<a:lin scaled="0"/>".
# NOTE: like classic VBA, .Color.RGB takes 0xBBGGRR (byte-swapped vs. the usual
# 0xRRGGBB hex notation used in the OOXML <a:srgbClr val="RRGGBB"/>), so the
# two target colors (F8F9FA / E8F5E8) are written byte-swapped below.
$bgShape.Fill.TwoColorGradient(2, 1)
$bgShape.Fill.GradientStops.Item(1).Color.RGB = 0xFAF9F8
$bgShape.Fill.GradientStops.Item(2).Color.RGB = 0xE8F5E8

# ---------------------------------------------------------------------------
# Helper: replace a whole paragraph's text with a single run. Writing the new
# text directly on top of the old often keeps a shared prefix/suffix as a
# separate, untouched run (PowerPoint's own "smart" diffing); routing through
# a disjoint placeholder first forces a full-text replace, collapsing the
# paragraph back down to exactly one <a:r>, matching authored output.
# ---------------------------------------------------------------------------
function Set-ParagraphText($textRange, [int]$paraIndex, [string]$newText) {
    $textRange.Paragraphs($paraIndex, 1).Text = "#"
    $textRange.Paragraphs($paraIndex, 1).Text = $newText
}

# ---------------------------------------------------------------------------
# 2) Slide 3 narrative text
# ---------------------------------------------------------------------------
$slide3Text = $slide3.Shapes.Item(4).TextFrame.TextRange

Set-ParagraphText $slide3Text 1 "За период 01.09.2025—07.09.2025 количество перезвонить по заявкам составило 111 из 164 (67,7%), что свидетельствует о высокой конверсии. Количество обработанных заявок превысило плановое на 12,5% (18 из 16). Однако объём одобренных гарантий (186 млн) значительно превышает выданный (39 млн)."
Set-ParagraphText $slide3Text 3 "Узкое место: низкий процент выдачи гарантий по отношению к одобренным."
Set-ParagraphText $slide3Text 6 "1. Проанализировать причины низкой выдачи по одобренным заявкам."
Set-ParagraphText $slide3Text 7 "2. Оптимизировать процесс выдачи гарантий для сокращения разрыва между одобрением и выдачей."
Set-ParagraphText $slide3Text 8 "3. Провести встречу с отделом обработки заявок для обсуждения текущих проблем и выработки решений."

# ---------------------------------------------------------------------------
# 3) Slide 9 narrative text
# ---------------------------------------------------------------------------
$slide9 = $p.Slides.Item(9)
$slide9Text = $slide9.Shapes.Item(4).TextFrame.TextRange

Set-ParagraphText $slide9Text 3 "За период с 01.09.2025 по 07.09.2025 количество перезвонить по заявкам составило 111 из 164 (67,7%), что является положительным результатом. Однако количество обработанных заявок превысило плановое (18 из 16, 112,5%) и сумма заявок значительно выше ожидаемой (646,0 из 204,0, 316,7%)."
Set-ParagraphText $slide9Text 5 "Узкое место: необходимо оптимизировать процесс обработки заявок для соответствия плановым показателям."
Set-ParagraphText $slide9Text 8 "1. Проанализировать причины превышения плановых показателей по заявкам."
Set-ParagraphText $slide9Text 9 "2. Скорректировать процессы обработки заявок для достижения плановых показателей."
Set-ParagraphText $slide9Text 10 "3. Провести анализ эффективности работы отдела и выявить возможности для оптимизации."

Write-Output "done"
